$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 12:52"

# --- Column A province/city name swaps (shared-string reorder in the
#     original file effectively swapped the text shown at these rows) ---
$ws.Range("A13").Value = "Alacant/Alicante"
$ws.Range("A14").Value = "Zaragoza"
$ws.Range("A15").Value = "Toledo"
$ws.Range("A16").Value = "Araba/Alava"

$ws.Range("A18").Value = "Valladolid"
$ws.Range("A19").Value = "Salamanca"
$ws.Range("A20").Value = "Malaga"

$ws.Range("A36").Value = "Castello/Castellon"
$ws.Range("A37").Value = "Guadalajara"
$ws.Range("A38").Value = "Soria"
$ws.Range("A39").Value = "Cadiz"
$ws.Range("A40").Value = "Avila"
$ws.Range("A41").Value = "Aragon"

# --- Numeric data updates (Casos totales / Casos activos / Recuperados / Muertes) ---
$data = @{
    8  = @(5267, 2365, 9790, 533)
    9  = @(4614, 1388, 2822, 404)
    11 = @(3450, 2365, 9790, 315)
    13 = @(3133, 990, 1812, 331)
    14 = @(3057, 706, 2017, 334)
    15 = @(3020, 2365, 9790, 414)
    16 = @(2990, 5092, 4658, 254)
    18 = @(2357, 833, 1320, 204)
    19 = @(2193, 607, 1334, 252)
    20 = @(2143, 523, 1450, 170)
    25 = @(1845, 508, 1191, 146)
    29 = @(1575, 781, 549, 245)
    33 = @(1207, 509, 558, 140)
    36 = @(1089, 246, 740, 103)
    37 = @(1077, 2365, 9790, 148)
    38 = @(990, 237, 670, 83)
    39 = @(943, 157, 734, 52)
    40 = @(917, 350, 473, 94)
    41 = @(907, 29, 838, 40)
    42 = @(884, 2365, 9790, 133)
    45 = @(620, 174, 398, 48)
    50 = @(414, 155, 208, 51)
    54 = @(98, 18, 78, 2)
    55 = @(93, 17, 72, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
